$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---------------------------------------------------
$meta = $wb.Worksheets.Item(1)

# Date: 2023-03-21T11:43:18+00:00 -> 2023-04-04T14:29:08+00:00
$meta.Range("B8").Value = "2023-04-04T14:29:08+00:00"

# Count: 2 -> 4 (kept as text, matching the rest of the "Value" column;
# a plain .Value assignment of "4" would be auto-coerced into a number,
# so build it as a text formula result first and paste just the value in).
$meta.Range("Z1").Formula = "=""4"""
$meta.Range("Z1").Copy()
$meta.Range("B21").PasteSpecial(-4163)  # xlPasteValues
$meta.Range("Z1").ClearContents()

# --- Concepts sheet -----------------------------------------------------
$concepts = $wb.Worksheets.Item(2)

# TARS -> TGS (row 4, column B "Code")
$concepts.Range("B4").Value = "TGS"

# Append a new row 5 for the "Whole Transcriptome Sequencing" concept,
# matching the formatting of the row above it (row 4).
$concepts.Range("A4:D4").Copy()
$concepts.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats

# Copy the "Level" value (and its underlying string type) from A4 into A5
# so that A5 keeps being a text cell with value "1", just like A2:A4.
$concepts.Range("A4").Copy()
$concepts.Range("A5").PasteSpecial(-4163)  # xlPasteValues

$concepts.Range("B5").Value = "WTS"
$concepts.Range("C5").Value = "Whole Transcriptome Sequencing"
# D5 (Definition) is intentionally left blank, same as the source row.

$excel.CutCopyMode = 0
